# Commit: "Encore des chgmts de titres"
# The "Year of Treatment" column (B) is removed entirely (all cells shift
# left by one column), and the remaining header cells (now B1:H1) get a
# ".jamais.jamais" suffix appended to their existing title text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the whole "Year of Treatment" column (column B). This shifts
#    every column C..I left by one, so former column I ("Total" / totals)
#    becomes column H, matching the new dimension A1:H33.
$ws.Columns.Item(2).Delete()

# 2) Append ".jamais.jamais" to every header cell from B1 to H1 (column A,
#    "Country", is left untouched).
for ($col = 2; $col -le 8; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value2()
    $cell.Value = "$current.jamais.jamais"
}
